$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.726.83'
$ws.Range('E2').Value = '  +1.76%  '

$ws.Range('D3').Value = '1.898.36'
$ws.Range('E3').Value = '  +2.56%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.20'
$ws.Range('E5').Value = '  +1.26%  '

$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4830'
$ws.Range('E7').Value = '  +1.00%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2855'
$ws.Range('E8').Value = '  +1.67%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06563'
$ws.Range('E9').Value = '  +1.26%  '

$ws.Range('D10').Value = '2.007.72'
$ws.Range('E10').Value = '  +8.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07465'
$ws.Range('E11').Value = '  +2.13%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.78'
$ws.Range('E12').Value = '  +3.01%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.118'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.16'
$ws.Range('E14').Value = '  +1.11%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6687'
$ws.Range('E15').Value = '  +3.40%  '

$ws.Range('D16').Value = '30.710.53'

$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '2.259.16'
$ws.Range('E17').Value = '  +7.74%  '

$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.34'
$ws.Range('E18').Value = '  +0.69%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.14%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007611'
$ws.Range('E20').Value = '  -0.15%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '231.61'
$ws.Range('E21').Value = '  +2.86%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.293'
$ws.Range('E22').Value = '  +0.11%  '

$ws.Range('E23').Value = '  +0.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.238'
$ws.Range('E24').Value = '  +2.63%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.06'
$ws.Range('E25').Value = '  +3.98%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.357'
$ws.Range('E26').Value = '  +1.45%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.81'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.969'
$ws.Range('E28').Value = '  +2.80%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.403'
$ws.Range('E29').Value = '  -1.89%  '

$ws.Range('E30').Value = '  +11.01%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.356'
$ws.Range('E31').Value = '  +2.60%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.033'
$ws.Range('E32').Value = '  +2.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05125'
$ws.Range('E33').Value = '  +2.25%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.221'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7604'
$ws.Range('E35').Value = '  +2.85%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.710'
$ws.Range('E36').Value = '  +0.93%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01885'
$ws.Range('E37').Value = '  +4.03%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.659'
$ws.Range('E38').Value = '  +1.86%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9211'
$ws.Range('E39').Value = '  +1.58%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.082'
$ws.Range('E40').Value = '  +1.46%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.38'
$ws.Range('E41').Value = '  +0.89%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4308'
$ws.Range('E42').Value = '  +1.27%  '

$ws.Range('E43').Value = '  +0.61%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.724'
$ws.Range('E44').Value = '  -3.86%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.442'
$ws.Range('E45').Value = '  +0.83%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.79'
$ws.Range('E46').Value = '  +1.11%  '

$ws.Range('E47').Value = '  -3.27%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.494'
$ws.Range('E48').Value = '  -3.64%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.964'
$ws.Range('E49').Value = '  +1.97%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.96'
$ws.Range('E50').Value = '  -0.75%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05679'
